# Refresh crypto Price (D) and Volume(1h) (E) columns with the latest feed values.
# Prices that are plain decimals (e.g. "1.007") get a leading apostrophe so Excel
# keeps storing them as text instead of silently coercing them to numbers -
# matching how the source feed always rendered this column as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = "23.841.98"; Volume = "  -3.30%  " },
    @{ Row = 3; Price = "1.622.51"; Volume = "  -3.22%  " },
    @{ Row = 4; Price = "1.007"; Volume = "  +0.53%  " },
    @{ Row = 5; Price = "1.005"; Volume = "  +0.28%  " },
    @{ Row = 6; Price = "306.31"; Volume = "  -2.55%  " },
    @{ Row = 7; Price = "0.3911"; Volume = "  +0.33%  " },
    @{ Row = 8; Price = "0.3820"; Volume = "  -3.06%  " },
    @{ Row = 9; Price = "1.005"; Volume = "  +0.37%  " },
    @{ Row = 10; Price = "49.76"; Volume = "  -4.07%  " },
    @{ Row = 11; Price = "1.358"; Volume = "  -2.42%  " },
    @{ Row = 12; Price = "0.08394"; Volume = "  -2.84%  " },
    @{ Row = 13; Price = "23.72"; Volume = "  -5.64%  " },
    @{ Row = 14; Price = "6.961"; Volume = "  -4.57%  " },
    @{ Row = 15; Price = "0.00001264"; Volume = "  -3.91%  " },
    @{ Row = 16; Price = "7.404"; Volume = "  -4.74%  " },
    @{ Row = 17; Price = "1.618.82"; Volume = "  -5.54%  " },
    @{ Row = 18; Price = "92.56"; Volume = "  -1.16%  " },
    @{ Row = 19; Price = "0.06906"; Volume = "  -2.17%  " },
    @{ Row = 20; Price = "19.79"; Volume = "  -3.95%  " },
    @{ Row = 21; Price = "6.814"; Volume = "  -3.49%  " },
    @{ Row = 22; Price = "1.006"; Volume = "  +0.19%  " },
    @{ Row = 23; Price = "13.32"; Volume = "  -5.07%  " },
    @{ Row = 24; Price = "23.863.80"; Volume = "  -3.23%  " },
    @{ Row = 25; Price = "2.385"; Volume = "  +0.91%  " },
    @{ Row = 26; Price = "2.853"; Volume = "  +4.66%  " },
    @{ Row = 27; Price = "22.01"; Volume = "  -4.94%  " },
    @{ Row = 28; Price = "157.70"; Volume = "  -2.70%  " },
    @{ Row = 29; Price = "138.65"; Volume = "  -5.45%  " },
    @{ Row = 30; Price = "5.226"; Volume = "  -9.15%  " },
    @{ Row = 31; Price = "7.637"; Volume = "  -3.42%  " },
    @{ Row = 32; Price = $null; Volume = "  -4.21%  " },
    @{ Row = 33; Price = "1.805.29"; Volume = "  -4.87%  " },
    @{ Row = 34; Price = "0.07905"; Volume = "  -5.44%  " },
    @{ Row = 35; Price = "0.02865"; Volume = "  -5.68%  " },
    @{ Row = 36; Price = "0.9515"; Volume = "  -3.37%  " },
    @{ Row = 37; Price = "6.539"; Volume = "  -4.44%  " },
    @{ Row = 38; Price = "0.2644"; Volume = "  -6.19%  " },
    @{ Row = 39; Price = "0.09106"; Volume = "  -4.10%  " },
    @{ Row = 40; Price = "10.29"; Volume = "  -1.86%  " },
    @{ Row = 41; Price = "13.42"; Volume = "  -0.78%  " },
    @{ Row = 42; Price = "1.415"; Volume = "  -8.14%  " },
    @{ Row = 43; Price = "0.7413"; Volume = "  -6.11%  " },
    @{ Row = 44; Price = "16.00"; Volume = "  -2.68%  " },
    @{ Row = 45; Price = "0.6797"; Volume = "  -4.58%  " },
    @{ Row = 46; Price = "2.428"; Volume = "  -5.08%  " },
    @{ Row = 47; Price = "4.053"; Volume = "  -3.25%  " },
    @{ Row = 48; Price = $null; Volume = "  +0.32%  " },
    @{ Row = 49; Price = "0.08228"; Volume = "  -4.82%  " },
    @{ Row = 50; Price = "132.17"; Volume = "  -3.70%  " },
    @{ Row = 51; Price = "1.243"; Volume = "  -6.01%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceText = $u.Price
        if ($priceText -match '^[+-]?\d+(\.\d+)?$') {
            # Prefix with an apostrophe so the numeric-looking text isn't
            # auto-converted into a Number by Excel.
            $priceText = "'" + $priceText
        }
        $ws.Range("D" + $u.Row).Value = $priceText
    }
    $ws.Range("E" + $u.Row).Value = $u.Volume
}
